$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 6.044933333333333
$ws.Range("H2").Value = 18.1348
$ws.Range("I2").Value = 0.9708761253868625
$ws.Range("J2").Value = 0.9708761253868624
$ws.Range("M2").Value = 2.874243666666667
$ws.Range("N2").Value = 8.622731
$ws.Range("O2").Value = 0.1078962025420678
$ws.Range("P2").Value = 0.1078962025420678
$ws.Range("Q2").Value = 17.37461134875555
$ws.Range("R2").Value = 156.3715021388
$ws.Range("S2").Value = 0.1047538470679989
$ws.Range("T2").Value = 0.1047538470679989

# Row 3
$ws.Range("G3").Value = 6.044933333333333
$ws.Range("H3").Value = 18.1348
$ws.Range("I3").Value = 0.9708761253868625
$ws.Range("J3").Value = 0.9708761253868624
$ws.Range("M3").Value = 5.803280666666666
$ws.Range("O3").Value = 0.2178492914434416
$ws.Range("P3").Value = 0.2178492914434416
$ws.Range("Q3").Value = 35.08044474462221
$ws.Range("R3").Value = 315.7240027015999
$ws.Range("S3").Value = 0.211504675994882
$ws.Range("T3").Value = 0.2115046759948819

# Row 4
$ws.Range("G4").Value = 6.044933333333333
$ws.Range("H4").Value = 18.1348
$ws.Range("I4").Value = 0.9708761253868625
$ws.Range("J4").Value = 0.9708761253868624
$ws.Range("M4").Value = 5.697719999999999
$ws.Range("N4").Value = 17.09316
$ws.Range("O4").Value = 0.2138866506961624
$ws.Range("P4").Value = 0.2138866506961624
$ws.Range("Q4").Value = 34.44233755199999
$ws.Range("R4").Value = 309.981037968
$ws.Range("S4").Value = 0.2076574426998635
$ws.Range("T4").Value = 0.2076574426998634

# Row 5
$ws.Range("G5").Value = 6.044933333333333
$ws.Range("H5").Value = 18.1348
$ws.Range("I5").Value = 0.9708761253868625
$ws.Range("J5").Value = 0.9708761253868624
$ws.Range("M5").Value = 1.156459333333333
$ws.Range("N5").Value = 3.469378
$ws.Range("O5").Value = 0.04341231465796556
$ws.Range("P5").Value = 0.04341231465796556
$ws.Range("Q5").Value = 6.99071957271111
$ws.Range("R5").Value = 62.91647615439999
$ws.Range("S5").Value = 0.04214797984920091
$ws.Range("T5").Value = 0.0421479798492009

# Row 6
$ws.Range("G6").Value = 6.044933333333333
$ws.Range("H6").Value = 18.1348
$ws.Range("I6").Value = 0.9708761253868625
$ws.Range("J6").Value = 0.9708761253868624
$ws.Range("M6").Value = 6.571656999999999
$ws.Range("N6").Value = 19.714971
$ws.Range("O6").Value = 0.24669336247727
$ws.Range("P6").Value = 0.24669336247727
$ws.Range("Q6").Value = 39.72522845453332
$ws.Range("R6").Value = 357.5270560907999
$ws.Range("S6").Value = 0.2395086959205887
$ws.Range("T6").Value = 0.2395086959205887

# Row 7
$ws.Range("G7").Value = 6.044933333333333
$ws.Range("H7").Value = 18.1348
$ws.Range("I7").Value = 0.9708761253868625
$ws.Range("J7").Value = 0.9708761253868624
$ws.Range("M7").Value = 4.535609
$ws.Range("N7").Value = 13.606827
$ws.Range("O7").Value = 0.1702621781830927
$ws.Range("P7").Value = 0.1702621781830927
$ws.Range("Q7").Value = 27.41745403106667
$ws.Range("R7").Value = 246.7570862796
$ws.Range("S7").Value = 0.1653034838543286
$ws.Range("T7").Value = 0.1653034838543286

# Row 8
$ws.Range("I8").Value = 0.01821359071319307
$ws.Range("J8").Value = 0.01821359071319307
$ws.Range("M8").Value = 2.874243666666667
$ws.Range("N8").Value = 8.622731
$ws.Range("O8").Value = 0.1078962025420678
$ws.Range("P8").Value = 0.1078962025420678
$ws.Range("Q8").Value = 0.3259468964497778
$ws.Range("R8").Value = 2.933522068048
$ws.Range("S8").Value = 0.001965177272609005
$ws.Range("T8").Value = 0.001965177272609005

# Row 9
$ws.Range("I9").Value = 0.01821359071319307
$ws.Range("J9").Value = 0.01821359071319307
$ws.Range("M9").Value = 5.803280666666666
$ws.Range("O9").Value = 0.2178492914434416
$ws.Range("P9").Value = 0.2178492914434416
$ws.Range("Q9").Value = 0.658107503015111
$ws.Range("S9").Value = 0.003967817831509958
$ws.Range("T9").Value = 0.003967817831509958

# Row 10
$ws.Range("I10").Value = 0.01821359071319307
$ws.Range("J10").Value = 0.01821359071319307
$ws.Range("M10").Value = 5.697719999999999
$ws.Range("N10").Value = 17.09316
$ws.Range("O10").Value = 0.2138866506961624
$ws.Range("P10").Value = 0.2138866506961624
$ws.Range("Q10").Value = 0.64613664192
$ws.Range("R10").Value = 5.815229777279999
$ws.Range("S10").Value = 0.003895643914795595
$ws.Range("T10").Value = 0.003895643914795595

# Row 11
$ws.Range("I11").Value = 0.01821359071319307
$ws.Range("J11").Value = 0.01821359071319307
$ws.Range("M11").Value = 1.156459333333333
$ws.Range("N11").Value = 3.469378
$ws.Range("O11").Value = 0.04341231465796556
$ws.Range("P11").Value = 0.04341231465796556
$ws.Range("Q11").Value = 0.1311455722915555
$ws.Range("R11").Value = 1.180310150624
$ws.Range("S11").Value = 0.0007906941310925371
$ws.Range("T11").Value = 0.0007906941310925371

# Row 12
$ws.Range("I12").Value = 0.01821359071319307
$ws.Range("J12").Value = 0.01821359071319307
$ws.Range("M12").Value = 6.571656999999999
$ws.Range("N12").Value = 19.714971
$ws.Range("O12").Value = 0.24669336247727
$ws.Range("P12").Value = 0.24669336247727
$ws.Range("Q12").Value = 0.7452434282186665
$ws.Range("R12").Value = 6.707190853968
$ws.Range("S12").Value = 0.004493171935822377
$ws.Range("T12").Value = 0.004493171935822377

# Row 13
$ws.Range("I13").Value = 0.01821359071319307
$ws.Range("J13").Value = 0.01821359071319307
$ws.Range("M13").Value = 4.535609
$ws.Range("N13").Value = 13.606827
$ws.Range("O13").Value = 0.1702621781830927
$ws.Range("P13").Value = 0.1702621781830927
$ws.Range("Q13").Value = 0.5143501555573333
$ws.Range("R13").Value = 4.629151400016
$ws.Range("S13").Value = 0.003101085627363601
$ws.Range("T13").Value = 0.003101085627363601

# Row 14
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.06793033333333333
$ws.Range("H14").Value = 0.203791
$ws.Range("I14").Value = 0.01091028389994453
$ws.Range("J14").Value = 0.01091028389994453
$ws.Range("M14").Value = 2.874243666666667
$ws.Range("N14").Value = 8.622731
$ws.Range("O14").Value = 0.1078962025420678
$ws.Range("P14").Value = 0.1078962025420678
$ws.Range("Q14").Value = 0.1952483303578889
$ws.Range("R14").Value = 1.757234973221
$ws.Range("S14").Value = 0.001177178201459876
$ws.Range("T14").Value = 0.001177178201459877

# Row 15
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.06793033333333333
$ws.Range("H15").Value = 0.203791
$ws.Range("I15").Value = 0.01091028389994453
$ws.Range("J15").Value = 0.01091028389994453
$ws.Range("M15").Value = 5.803280666666666
$ws.Range("O15").Value = 0.2178492914434416
$ws.Range("P15").Value = 0.2178492914434416
$ws.Range("Q15").Value = 0.3942187901135555
$ws.Range("R15").Value = 3.547969111021999
$ws.Range("S15").Value = 0.002376797617049705
$ws.Range("T15").Value = 0.002376797617049705

# Row 16
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.06793033333333333
$ws.Range("H16").Value = 0.203791
$ws.Range("I16").Value = 0.01091028389994453
$ws.Range("J16").Value = 0.01091028389994453
$ws.Range("M16").Value = 5.697719999999999
$ws.Range("N16").Value = 17.09316
$ws.Range("O16").Value = 0.2138866506961624
$ws.Range("P16").Value = 0.2138866506961624
$ws.Range("Q16").Value = 0.38704801884
$ws.Range("R16").Value = 3.483432169559999
$ws.Range("S16").Value = 0.002333564081503401
$ws.Range("T16").Value = 0.002333564081503401

# Row 17
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.06793033333333333
$ws.Range("H17").Value = 0.203791
$ws.Range("I17").Value = 0.01091028389994453
$ws.Range("J17").Value = 0.01091028389994453
$ws.Range("M17").Value = 1.156459333333333
$ws.Range("N17").Value = 3.469378
$ws.Range("O17").Value = 0.04341231465796556
$ws.Range("P17").Value = 0.04341231465796556
$ws.Range("Q17").Value = 0.07855866799977777
$ws.Range("R17").Value = 0.7070280119979999
$ws.Range("S17").Value = 0.0004736406776721277
$ws.Range("T17").Value = 0.0004736406776721277

# Row 18
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.06793033333333333
$ws.Range("H18").Value = 0.203791
$ws.Range("I18").Value = 0.01091028389994453
$ws.Range("J18").Value = 0.01091028389994453
$ws.Range("M18").Value = 6.571656999999999
$ws.Range("N18").Value = 19.714971
$ws.Range("O18").Value = 0.24669336247727
$ws.Range("P18").Value = 0.24669336247727
$ws.Range("Q18").Value = 0.4464148505623333
$ws.Range("R18").Value = 4.017733655061
$ws.Range("S18").Value = 0.002691494620858939
$ws.Range("T18").Value = 0.002691494620858939

# Row 19
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.06793033333333333
$ws.Range("H19").Value = 0.203791
$ws.Range("I19").Value = 0.01091028389994453
$ws.Range("J19").Value = 0.01091028389994453
$ws.Range("M19").Value = 4.535609
$ws.Range("N19").Value = 13.606827
$ws.Range("O19").Value = 0.1702621781830927
$ws.Range("P19").Value = 0.1702621781830927
$ws.Range("Q19").Value = 0.3081054312396667
$ws.Range("R19").Value = 2.772948881157
$ws.Range("S19").Value = 0.1857608701400483
$ws.Range("T19").Value = 0.1857608701400483
